$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 (extr1 .. extr8) with new C/D/E values ---
# Row 8 (extr1): from_bus 5->14, to_bus 12->11, in_service FALSE->TRUE
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9 (extr2): from_bus 5->16 (to_bus/in_service unchanged)
$ws.Range("C9").Value = 16

# Row 10 (extr3): from_bus 10->5, to_bus 11->12 (in_service unchanged)
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

# Row 11 (extr4): from_bus 7->5, to_bus 8->9 (in_service unchanged)
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

# Row 12 (extr5): from_bus 9->10, in_service FALSE->TRUE (to_bus unchanged)
$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $true

# Row 13 (extr6): to_bus 11->8, in_service FALSE->TRUE (from_bus unchanged)
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# Row 14 (extr7): from_bus 5->9, to_bus 7->11, in_service TRUE->FALSE
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# Row 15 (extr8): from_bus 8->7, to_bus 5->11, in_service TRUE->FALSE
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- Append two new rows (line7, line8) copying the formatting of row 15 ---
$ws.Range("A15:E15").Copy($ws.Range("A16:E16"))
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A15:E15").Copy($ws.Range("A17:E17"))
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
